# Fill in the cfu_count_undiluted values (column C) for rows 2-33,
# then move the active selection to F2 (matching the author's final click).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 5500000
    3  = 7200000
    4  = 5300000
    5  = 3500000
    6  = 38000000
    7  = 39000000
    8  = 35000000
    9  = 43000000
    10 = 13000000
    11 = 17000000
    12 = 12000000
    13 = 16000000
    14 = 690000
    15 = 770000
    16 = 710000
    17 = 810000
    18 = 470000
    19 = 520000
    20 = 560000
    21 = 420000
    22 = 410000
    23 = 510000
    24 = 430000
    25 = 410000
    26 = 4600000
    27 = 2800000
    28 = 3500000
    29 = 4000000
    30 = 3700000
    31 = 5300000
    32 = 5200000
    33 = 4000000
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}

$ws.Range("F2").Select()
